$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 18 (shifts existing rows 18-68 down to 19-69)
$ws.Rows.Item(18).Insert()

# Populate the new row 18 with the new weekly price record
$ws.Range("A18").Value = 7
$ws.Range("B18").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C18").Value = "Ñuble"
$ws.Range("D18").Value = 44620
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 100112022
$ws.Range("G18").Value = "Arveja Verde"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 24000
$ws.Range("L18").Value = 25000
$ws.Range("M18").Value = 24500
$ws.Range("N18").Value = "`$/saco 25 kilos"
$ws.Range("O18").Value = "Provincia de Diguillín"
$ws.Range("P18").Value = 980
$ws.Range("Q18").Value = 25
$ws.Range("R18").Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Range("D18").NumberFormat = $ws.Range("D19").NumberFormat
